$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new columns before column D (shifts old D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Step 2: Copy number formatting from column F into the newly inserted D:E columns
# so the blank cells created by Insert() pick up the correct style (date/number format)
# for every row that has data in the source sheet (rows 5 to 102).
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: Write the final values into D:M for every data row.
# Row 7
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(7, 6).Value = 43281
$ws.Cells.Item(7, 7).Value = 43190
$ws.Cells.Item(7, 8).Value = 43100
$ws.Cells.Item(7, 9).Value = 43008
$ws.Cells.Item(7, 10).Value = 42916
$ws.Cells.Item(7, 11).Value = 42825
$ws.Cells.Item(7, 12).Value = 42735
$ws.Cells.Item(7, 13).Value = 42643

# Row 8
$ws.Cells.Item(8, 4).Value = 689100
$ws.Cells.Item(8, 5).Value = 674800
$ws.Cells.Item(8, 6).Value = 1332100
$ws.Cells.Item(8, 7).Value = 659400
$ws.Cells.Item(8, 8).Value = 675800
$ws.Cells.Item(8, 9).Value = 657200
$ws.Cells.Item(8, 10).Value = 1258800
$ws.Cells.Item(8, 11).Value = 582700
$ws.Cells.Item(8, 12).Value = 629600
$ws.Cells.Item(8, 13).Value = 621600

# Row 9
$ws.Cells.Item(9, 4).Value = 432700
$ws.Cells.Item(9, 5).Value = 448900
$ws.Cells.Item(9, 6).Value = 805800
$ws.Cells.Item(9, 7).Value = 365200
$ws.Cells.Item(9, 8).Value = 371900
$ws.Cells.Item(9, 9).Value = 375000
$ws.Cells.Item(9, 10).Value = 707400
$ws.Cells.Item(9, 11).Value = 325700
$ws.Cells.Item(9, 12).Value = 300100
$ws.Cells.Item(9, 13).Value = 293200

# Row 10
$ws.Cells.Item(10, 4).Value = 256400
$ws.Cells.Item(10, 5).Value = 225900
$ws.Cells.Item(10, 6).Value = 526300
$ws.Cells.Item(10, 7).Value = 294200
$ws.Cells.Item(10, 8).Value = 303900
$ws.Cells.Item(10, 9).Value = 282200
$ws.Cells.Item(10, 10).Value = 551400
$ws.Cells.Item(10, 11).Value = 257000
$ws.Cells.Item(10, 12).Value = 329500
$ws.Cells.Item(10, 13).Value = 328400

# Row 11
$ws.Cells.Item(11, 4).Value = ""
$ws.Cells.Item(11, 5).Value = ""
$ws.Cells.Item(11, 6).Value = ""
$ws.Cells.Item(11, 7).Value = ""
$ws.Cells.Item(11, 8).Value = ""
$ws.Cells.Item(11, 9).Value = ""
$ws.Cells.Item(11, 10).Value = ""
$ws.Cells.Item(11, 11).Value = ""
$ws.Cells.Item(11, 12).Value = ""
$ws.Cells.Item(11, 13).Value = ""

# Row 12
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(12, 6).Value = "NA"
$ws.Cells.Item(12, 7).Value = "NA"
$ws.Cells.Item(12, 8).Value = "NA"
$ws.Cells.Item(12, 9).Value = "NA"
$ws.Cells.Item(12, 10).Value = "NA"
$ws.Cells.Item(12, 11).Value = "NA"
$ws.Cells.Item(12, 12).Value = "NA"
$ws.Cells.Item(12, 13).Value = "NA"

# Row 13
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0

# Row 14
$ws.Cells.Item(14, 4).Value = 2700
$ws.Cells.Item(14, 5).Value = 7200
$ws.Cells.Item(14, 6).Value = 4100
$ws.Cells.Item(14, 7).Value = -5900
$ws.Cells.Item(14, 8).Value = 500
$ws.Cells.Item(14, 9).Value = 300
$ws.Cells.Item(14, 10).Value = 9600
$ws.Cells.Item(14, 11).Value = 300
$ws.Cells.Item(14, 12).Value = 1600
$ws.Cells.Item(14, 13).Value = 0

# Row 15
$ws.Cells.Item(15, 4).Value = 28200
$ws.Cells.Item(15, 5).Value = 25900
$ws.Cells.Item(15, 6).Value = 50000
$ws.Cells.Item(15, 7).Value = 25000
$ws.Cells.Item(15, 8).Value = 28900
$ws.Cells.Item(15, 9).Value = 26900
$ws.Cells.Item(15, 10).Value = 47900
$ws.Cells.Item(15, 11).Value = 21500
$ws.Cells.Item(15, 12).Value = 23200
$ws.Cells.Item(15, 13).Value = 22700

# Row 16
$ws.Cells.Item(16, 4).Value = ""
$ws.Cells.Item(16, 5).Value = ""
$ws.Cells.Item(16, 6).Value = ""
$ws.Cells.Item(16, 7).Value = ""
$ws.Cells.Item(16, 8).Value = ""
$ws.Cells.Item(16, 9).Value = ""
$ws.Cells.Item(16, 10).Value = ""
$ws.Cells.Item(16, 11).Value = ""
$ws.Cells.Item(16, 12).Value = ""
$ws.Cells.Item(16, 13).Value = ""

# Row 17
$ws.Cells.Item(17, 4).Value = 616200
$ws.Cells.Item(17, 5).Value = 613100
$ws.Cells.Item(17, 6).Value = 1226400
$ws.Cells.Item(17, 7).Value = 609300
$ws.Cells.Item(17, 8).Value = 626800
$ws.Cells.Item(17, 9).Value = 630300
$ws.Cells.Item(17, 10).Value = 1199200
$ws.Cells.Item(17, 11).Value = 572800
$ws.Cells.Item(17, 12).Value = 520100
$ws.Cells.Item(17, 13).Value = 553600

# Row 18
$ws.Cells.Item(18, 4).Value = 72900
$ws.Cells.Item(18, 5).Value = 61700
$ws.Cells.Item(18, 6).Value = 105700
$ws.Cells.Item(18, 7).Value = 50100
$ws.Cells.Item(18, 8).Value = 49000
$ws.Cells.Item(18, 9).Value = 26900
$ws.Cells.Item(18, 10).Value = 59600
$ws.Cells.Item(18, 11).Value = 9900
$ws.Cells.Item(18, 12).Value = 109500
$ws.Cells.Item(18, 13).Value = 68000

# Row 19
$ws.Cells.Item(19, 4).Value = ""
$ws.Cells.Item(19, 5).Value = ""
$ws.Cells.Item(19, 6).Value = ""
$ws.Cells.Item(19, 7).Value = ""
$ws.Cells.Item(19, 8).Value = ""
$ws.Cells.Item(19, 9).Value = ""
$ws.Cells.Item(19, 10).Value = ""
$ws.Cells.Item(19, 11).Value = ""
$ws.Cells.Item(19, 12).Value = ""
$ws.Cells.Item(19, 13).Value = ""

# Row 20
$ws.Cells.Item(20, 4).Value = 2800
$ws.Cells.Item(20, 5).Value = 79500
$ws.Cells.Item(20, 6).Value = 38800
$ws.Cells.Item(20, 7).Value = 14600
$ws.Cells.Item(20, 8).Value = 13900
$ws.Cells.Item(20, 9).Value = 19900
$ws.Cells.Item(20, 10).Value = 47200
$ws.Cells.Item(20, 11).Value = 22000
$ws.Cells.Item(20, 12).Value = -34500
$ws.Cells.Item(20, 13).Value = -18500

# Row 21
$ws.Cells.Item(21, 4).Value = 103800
$ws.Cells.Item(21, 5).Value = 175200
$ws.Cells.Item(21, 6).Value = 194600
$ws.Cells.Item(21, 7).Value = 89700
$ws.Cells.Item(21, 8).Value = 91900
$ws.Cells.Item(21, 9).Value = 74100
$ws.Cells.Item(21, 10).Value = 164000
$ws.Cells.Item(21, 11).Value = 53400
$ws.Cells.Item(21, 12).Value = 99800
$ws.Cells.Item(21, 13).Value = 72300

# Row 22
$ws.Cells.Item(22, 4).Value = 6500
$ws.Cells.Item(22, 5).Value = 6100
$ws.Cells.Item(22, 6).Value = 25200
$ws.Cells.Item(22, 7).Value = 8100
$ws.Cells.Item(22, 8).Value = 8100
$ws.Cells.Item(22, 9).Value = 8600
$ws.Cells.Item(22, 10).Value = 17200
$ws.Cells.Item(22, 11).Value = 8100
$ws.Cells.Item(22, 12).Value = 10900
$ws.Cells.Item(22, 13).Value = 8600

# Row 23
$ws.Cells.Item(23, 4).Value = 69100
$ws.Cells.Item(23, 5).Value = 135100
$ws.Cells.Item(23, 6).Value = 119300
$ws.Cells.Item(23, 7).Value = 56600
$ws.Cells.Item(23, 8).Value = 54900
$ws.Cells.Item(23, 9).Value = 38200
$ws.Cells.Item(23, 10).Value = 89700
$ws.Cells.Item(23, 11).Value = 23800
$ws.Cells.Item(23, 12).Value = 64100
$ws.Cells.Item(23, 13).Value = 40900

# Row 24
$ws.Cells.Item(24, 4).Value = 12400
$ws.Cells.Item(24, 5).Value = 10000
$ws.Cells.Item(24, 6).Value = 29700
$ws.Cells.Item(24, 7).Value = 13600
$ws.Cells.Item(24, 8).Value = -131400
$ws.Cells.Item(24, 9).Value = 13400
$ws.Cells.Item(24, 10).Value = 26600
$ws.Cells.Item(24, 11).Value = 2700
$ws.Cells.Item(24, 12).Value = 27200
$ws.Cells.Item(24, 13).Value = 7800

# Row 25
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0

# Row 26
$ws.Cells.Item(26, 4).Value = 56700
$ws.Cells.Item(26, 5).Value = 125100
$ws.Cells.Item(26, 6).Value = 89600
$ws.Cells.Item(26, 7).Value = 43000
$ws.Cells.Item(26, 8).Value = 186300
$ws.Cells.Item(26, 9).Value = 24800
$ws.Cells.Item(26, 10).Value = 63100
$ws.Cells.Item(26, 11).Value = 21100
$ws.Cells.Item(26, 12).Value = 36900
$ws.Cells.Item(26, 13).Value = 33100

# Row 27
$ws.Cells.Item(27, 4).Value = 56400
$ws.Cells.Item(27, 5).Value = 124200
$ws.Cells.Item(27, 6).Value = 88900
$ws.Cells.Item(27, 7).Value = 42600
$ws.Cells.Item(27, 8).Value = 183900
$ws.Cells.Item(27, 9).Value = 24600
$ws.Cells.Item(27, 10).Value = 62500
$ws.Cells.Item(27, 11).Value = 20900
$ws.Cells.Item(27, 12).Value = 36600
$ws.Cells.Item(27, 13).Value = 32700

# Row 28
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0

# Row 29
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(29, 6).Value = "NA"
$ws.Cells.Item(29, 7).Value = "NA"
$ws.Cells.Item(29, 8).Value = 28300
$ws.Cells.Item(29, 9).Value = "NA"
$ws.Cells.Item(29, 10).Value = "NA"
$ws.Cells.Item(29, 11).Value = "NA"
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0

# Row 30
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0

# Row 31
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 0

# Row 32
$ws.Cells.Item(32, 4).Value = -2800
$ws.Cells.Item(32, 5).Value = -79500
$ws.Cells.Item(32, 6).Value = -38800
$ws.Cells.Item(32, 7).Value = -14600
$ws.Cells.Item(32, 8).Value = -13900
$ws.Cells.Item(32, 9).Value = -19900
$ws.Cells.Item(32, 10).Value = -47200
$ws.Cells.Item(32, 11).Value = -22000
$ws.Cells.Item(32, 12).Value = 34500
$ws.Cells.Item(32, 13).Value = 18500

# Row 33
$ws.Cells.Item(33, 4).Value = 56400
$ws.Cells.Item(33, 5).Value = 124200
$ws.Cells.Item(33, 6).Value = 88900
$ws.Cells.Item(33, 7).Value = 42600
$ws.Cells.Item(33, 8).Value = 212200
$ws.Cells.Item(33, 9).Value = 24600
$ws.Cells.Item(33, 10).Value = 62500
$ws.Cells.Item(33, 11).Value = 20900
$ws.Cells.Item(33, 12).Value = 36600
$ws.Cells.Item(33, 13).Value = 32700

# Row 34
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 0

# Row 35
$ws.Cells.Item(35, 4).Value = 56400
$ws.Cells.Item(35, 5).Value = 124200
$ws.Cells.Item(35, 6).Value = 88900
$ws.Cells.Item(35, 7).Value = 42600
$ws.Cells.Item(35, 8).Value = 212200
$ws.Cells.Item(35, 9).Value = 24600
$ws.Cells.Item(35, 10).Value = 62500
$ws.Cells.Item(35, 11).Value = 20900
$ws.Cells.Item(35, 12).Value = 36600
$ws.Cells.Item(35, 13).Value = 32700

# Row 38
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(38, 6).Value = 43281
$ws.Cells.Item(38, 7).Value = 43190
$ws.Cells.Item(38, 8).Value = 43100
$ws.Cells.Item(38, 9).Value = 43008
$ws.Cells.Item(38, 10).Value = 42916
$ws.Cells.Item(38, 11).Value = 42825
$ws.Cells.Item(38, 12).Value = 42735
$ws.Cells.Item(38, 13).Value = 42643

# Row 39
$ws.Cells.Item(39, 4).Value = ""
$ws.Cells.Item(39, 5).Value = ""
$ws.Cells.Item(39, 6).Value = ""
$ws.Cells.Item(39, 7).Value = ""
$ws.Cells.Item(39, 8).Value = ""
$ws.Cells.Item(39, 9).Value = ""
$ws.Cells.Item(39, 10).Value = ""
$ws.Cells.Item(39, 11).Value = ""
$ws.Cells.Item(39, 12).Value = ""
$ws.Cells.Item(39, 13).Value = ""

# Row 40
$ws.Cells.Item(40, 4).Value = ""
$ws.Cells.Item(40, 5).Value = ""
$ws.Cells.Item(40, 6).Value = ""
$ws.Cells.Item(40, 7).Value = ""
$ws.Cells.Item(40, 8).Value = ""
$ws.Cells.Item(40, 9).Value = ""
$ws.Cells.Item(40, 10).Value = ""
$ws.Cells.Item(40, 11).Value = ""
$ws.Cells.Item(40, 12).Value = ""
$ws.Cells.Item(40, 13).Value = ""

# Row 41
$ws.Cells.Item(41, 4).Value = 253300
$ws.Cells.Item(41, 5).Value = 265200
$ws.Cells.Item(41, 6).Value = 288600
$ws.Cells.Item(41, 7).Value = 315800
$ws.Cells.Item(41, 8).Value = 390000
$ws.Cells.Item(41, 9).Value = 395000
$ws.Cells.Item(41, 10).Value = 432100
$ws.Cells.Item(41, 11).Value = 632800
$ws.Cells.Item(41, 12).Value = 648900
$ws.Cells.Item(41, 13).Value = 593200

# Row 42
$ws.Cells.Item(42, 4).Value = 514600
$ws.Cells.Item(42, 5).Value = 517900
$ws.Cells.Item(42, 6).Value = 474100
$ws.Cells.Item(42, 7).Value = 491400
$ws.Cells.Item(42, 8).Value = 557200
$ws.Cells.Item(42, 9).Value = 519600
$ws.Cells.Item(42, 10).Value = 470000
$ws.Cells.Item(42, 11).Value = 456600
$ws.Cells.Item(42, 12).Value = 448200
$ws.Cells.Item(42, 13).Value = 408700

# Row 43
$ws.Cells.Item(43, 4).Value = 601400
$ws.Cells.Item(43, 5).Value = 589600
$ws.Cells.Item(43, 6).Value = 563600
$ws.Cells.Item(43, 7).Value = 545900
$ws.Cells.Item(43, 8).Value = 644200
$ws.Cells.Item(43, 9).Value = 553900
$ws.Cells.Item(43, 10).Value = 545500
$ws.Cells.Item(43, 11).Value = 541500
$ws.Cells.Item(43, 12).Value = 656700
$ws.Cells.Item(43, 13).Value = 587300

# Row 44
$ws.Cells.Item(44, 4).Value = 69500
$ws.Cells.Item(44, 5).Value = 73700
$ws.Cells.Item(44, 6).Value = 76000
$ws.Cells.Item(44, 7).Value = 66900
$ws.Cells.Item(44, 8).Value = 60600
$ws.Cells.Item(44, 9).Value = 62500
$ws.Cells.Item(44, 10).Value = 58800
$ws.Cells.Item(44, 11).Value = 35700
$ws.Cells.Item(44, 12).Value = 34800
$ws.Cells.Item(44, 13).Value = 35300

# Row 45
$ws.Cells.Item(45, 4).Value = 93600
$ws.Cells.Item(45, 5).Value = 89400
$ws.Cells.Item(45, 6).Value = 92900
$ws.Cells.Item(45, 7).Value = 101400
$ws.Cells.Item(45, 8).Value = 83800
$ws.Cells.Item(45, 9).Value = 86800
$ws.Cells.Item(45, 10).Value = 91000
$ws.Cells.Item(45, 11).Value = 93300
$ws.Cells.Item(45, 12).Value = 82700
$ws.Cells.Item(45, 13).Value = 91200

# Row 46
$ws.Cells.Item(46, 4).Value = 1532300
$ws.Cells.Item(46, 5).Value = 1535900
$ws.Cells.Item(46, 6).Value = 1495200
$ws.Cells.Item(46, 7).Value = 1521400
$ws.Cells.Item(46, 8).Value = 1735800
$ws.Cells.Item(46, 9).Value = 1617900
$ws.Cells.Item(46, 10).Value = 1597400
$ws.Cells.Item(46, 11).Value = 1759900
$ws.Cells.Item(46, 12).Value = 1871300
$ws.Cells.Item(46, 13).Value = 1715600

# Row 47
$ws.Cells.Item(47, 4).Value = 143800
$ws.Cells.Item(47, 5).Value = 142800
$ws.Cells.Item(47, 6).Value = 133600
$ws.Cells.Item(47, 7).Value = 131900
$ws.Cells.Item(47, 8).Value = 128600
$ws.Cells.Item(47, 9).Value = 122200
$ws.Cells.Item(47, 10).Value = 67800
$ws.Cells.Item(47, 11).Value = 60300
$ws.Cells.Item(47, 12).Value = 58800
$ws.Cells.Item(47, 13).Value = 64800

# Row 48
$ws.Cells.Item(48, 4).Value = 293100
$ws.Cells.Item(48, 5).Value = 271600
$ws.Cells.Item(48, 6).Value = 265500
$ws.Cells.Item(48, 7).Value = 258700
$ws.Cells.Item(48, 8).Value = 259400
$ws.Cells.Item(48, 9).Value = 259800
$ws.Cells.Item(48, 10).Value = 263100
$ws.Cells.Item(48, 11).Value = 243300
$ws.Cells.Item(48, 12).Value = 233700
$ws.Cells.Item(48, 13).Value = 225700

# Row 49
$ws.Cells.Item(49, 4).Value = 1660000
$ws.Cells.Item(49, 5).Value = 1684900
$ws.Cells.Item(49, 6).Value = 1617400
$ws.Cells.Item(49, 7).Value = 1636100
$ws.Cells.Item(49, 8).Value = 1639900
$ws.Cells.Item(49, 9).Value = 1648300
$ws.Cells.Item(49, 10).Value = 1638300
$ws.Cells.Item(49, 11).Value = 1426600
$ws.Cells.Item(49, 12).Value = 1296900
$ws.Cells.Item(49, 13).Value = 1333900

# Row 50
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 0

# Row 51
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = 0

# Row 52
$ws.Cells.Item(52, 4).Value = 1134800
$ws.Cells.Item(52, 5).Value = 1238700
$ws.Cells.Item(52, 6).Value = 1219300
$ws.Cells.Item(52, 7).Value = 1212900
$ws.Cells.Item(52, 8).Value = 1174200
$ws.Cells.Item(52, 9).Value = 971800
$ws.Cells.Item(52, 10).Value = 941800
$ws.Cells.Item(52, 11).Value = 926000
$ws.Cells.Item(52, 12).Value = 971900
$ws.Cells.Item(52, 13).Value = 1098700

# Row 53
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 0

# Row 54
$ws.Cells.Item(54, 4).Value = 4764000
$ws.Cells.Item(54, 5).Value = 4873700
$ws.Cells.Item(54, 6).Value = 4730900
$ws.Cells.Item(54, 7).Value = 4761000
$ws.Cells.Item(54, 8).Value = 4937800
$ws.Cells.Item(54, 9).Value = 4619900
$ws.Cells.Item(54, 10).Value = 4508400
$ws.Cells.Item(54, 11).Value = 4416100
$ws.Cells.Item(54, 12).Value = 4432700
$ws.Cells.Item(54, 13).Value = 4438600

# Row 55
$ws.Cells.Item(55, 4).Value = ""
$ws.Cells.Item(55, 5).Value = ""
$ws.Cells.Item(55, 6).Value = ""
$ws.Cells.Item(55, 7).Value = ""
$ws.Cells.Item(55, 8).Value = ""
$ws.Cells.Item(55, 9).Value = ""
$ws.Cells.Item(55, 10).Value = ""
$ws.Cells.Item(55, 11).Value = ""
$ws.Cells.Item(55, 12).Value = ""
$ws.Cells.Item(55, 13).Value = ""

# Row 56
$ws.Cells.Item(56, 4).Value = ""
$ws.Cells.Item(56, 5).Value = ""
$ws.Cells.Item(56, 6).Value = ""
$ws.Cells.Item(56, 7).Value = ""
$ws.Cells.Item(56, 8).Value = ""
$ws.Cells.Item(56, 9).Value = ""
$ws.Cells.Item(56, 10).Value = ""
$ws.Cells.Item(56, 11).Value = ""
$ws.Cells.Item(56, 12).Value = ""
$ws.Cells.Item(56, 13).Value = ""

# Row 57
$ws.Cells.Item(57, 4).Value = 486600
$ws.Cells.Item(57, 5).Value = 448200
$ws.Cells.Item(57, 6).Value = 457300
$ws.Cells.Item(57, 7).Value = 435500
$ws.Cells.Item(57, 8).Value = 526300
$ws.Cells.Item(57, 9).Value = 446100
$ws.Cells.Item(57, 10).Value = 451500
$ws.Cells.Item(57, 11).Value = 444700
$ws.Cells.Item(57, 12).Value = 352400
$ws.Cells.Item(57, 13).Value = 406200

# Row 58
$ws.Cells.Item(58, 4).Value = 6400
$ws.Cells.Item(58, 5).Value = 6600
$ws.Cells.Item(58, 6).Value = 23100
$ws.Cells.Item(58, 7).Value = 406700
$ws.Cells.Item(58, 8).Value = 6700
$ws.Cells.Item(58, 9).Value = 6700
$ws.Cells.Item(58, 10).Value = 6500
$ws.Cells.Item(58, 11).Value = 6200
$ws.Cells.Item(58, 12).Value = 6100
$ws.Cells.Item(58, 13).Value = 6500

# Row 59
$ws.Cells.Item(59, 4).Value = 319200
$ws.Cells.Item(59, 5).Value = 349400
$ws.Cells.Item(59, 6).Value = 266000
$ws.Cells.Item(59, 7).Value = 310400
$ws.Cells.Item(59, 8).Value = 345600
$ws.Cells.Item(59, 9).Value = 360400
$ws.Cells.Item(59, 10).Value = 297300
$ws.Cells.Item(59, 11).Value = 328200
$ws.Cells.Item(59, 12).Value = 460500
$ws.Cells.Item(59, 13).Value = 341400

# Row 60
$ws.Cells.Item(60, 4).Value = 812200
$ws.Cells.Item(60, 5).Value = 804200
$ws.Cells.Item(60, 6).Value = 746400
$ws.Cells.Item(60, 7).Value = 1152600
$ws.Cells.Item(60, 8).Value = 878600
$ws.Cells.Item(60, 9).Value = 813200
$ws.Cells.Item(60, 10).Value = 755200
$ws.Cells.Item(60, 11).Value = 779000
$ws.Cells.Item(60, 12).Value = 819000
$ws.Cells.Item(60, 13).Value = 754100

# Row 61
$ws.Cells.Item(61, 4).Value = 470800
$ws.Cells.Item(61, 5).Value = 473100
$ws.Cells.Item(61, 6).Value = 478800
$ws.Cells.Item(61, 7).Value = 101400
$ws.Cells.Item(61, 8).Value = 496900
$ws.Cells.Item(61, 9).Value = 498800
$ws.Cells.Item(61, 10).Value = 502300
$ws.Cells.Item(61, 11).Value = 499800
$ws.Cells.Item(61, 12).Value = 498300
$ws.Cells.Item(61, 13).Value = 500800

# Row 62
$ws.Cells.Item(62, 4).Value = 560000
$ws.Cells.Item(62, 5).Value = 616400
$ws.Cells.Item(62, 6).Value = 623600
$ws.Cells.Item(62, 7).Value = 617400
$ws.Cells.Item(62, 8).Value = 642600
$ws.Cells.Item(62, 9).Value = 743500
$ws.Cells.Item(62, 10).Value = 710800
$ws.Cells.Item(62, 11).Value = 656700
$ws.Cells.Item(62, 12).Value = 662400
$ws.Cells.Item(62, 13).Value = 675000

# Row 63
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0

# Row 64
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = 0

# Row 65
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = 0

# Row 66
$ws.Cells.Item(66, 4).Value = 1847300
$ws.Cells.Item(66, 5).Value = 1898400
$ws.Cells.Item(66, 6).Value = 1853500
$ws.Cells.Item(66, 7).Value = 1876000
$ws.Cells.Item(66, 8).Value = 2022700
$ws.Cells.Item(66, 9).Value = 2059200
$ws.Cells.Item(66, 10).Value = 1972100
$ws.Cells.Item(66, 11).Value = 1935500
$ws.Cells.Item(66, 12).Value = 1979700
$ws.Cells.Item(66, 13).Value = 1929800

# Row 67
$ws.Cells.Item(67, 4).Value = ""
$ws.Cells.Item(67, 5).Value = ""
$ws.Cells.Item(67, 6).Value = ""
$ws.Cells.Item(67, 7).Value = ""
$ws.Cells.Item(67, 8).Value = ""
$ws.Cells.Item(67, 9).Value = ""
$ws.Cells.Item(67, 10).Value = ""
$ws.Cells.Item(67, 11).Value = ""
$ws.Cells.Item(67, 12).Value = ""
$ws.Cells.Item(67, 13).Value = ""

# Row 68
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 0

# Row 69
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = 0

# Row 70
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0

# Row 71
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 0

# Row 72
$ws.Cells.Item(72, 4).Value = 6236100
$ws.Cells.Item(72, 5).Value = 6179900
$ws.Cells.Item(72, 6).Value = 6061900
$ws.Cells.Item(72, 7).Value = 6022300
$ws.Cells.Item(72, 8).Value = 5791700
$ws.Cells.Item(72, 9).Value = 5648500
$ws.Cells.Item(72, 10).Value = 5630700
$ws.Cells.Item(72, 11).Value = 5595800
$ws.Cells.Item(72, 12).Value = 5588900
$ws.Cells.Item(72, 13).Value = 5552000

# Row 73
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = 0

# Row 74
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = 0

# Row 75
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).Value = 0

# Row 76
$ws.Cells.Item(76, 4).Value = 2916800
$ws.Cells.Item(76, 5).Value = 2975400
$ws.Cells.Item(76, 6).Value = 2877400
$ws.Cells.Item(76, 7).Value = 2884900
$ws.Cells.Item(76, 8).Value = 2915100
$ws.Cells.Item(76, 9).Value = 2560700
$ws.Cells.Item(76, 10).Value = 2536400
$ws.Cells.Item(76, 11).Value = 2480600
$ws.Cells.Item(76, 12).Value = 2452900
$ws.Cells.Item(76, 13).Value = 2508800

# Row 77
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = 0

# Row 80
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(80, 6).Value = 43281
$ws.Cells.Item(80, 7).Value = 43190
$ws.Cells.Item(80, 8).Value = 43100
$ws.Cells.Item(80, 9).Value = 43008
$ws.Cells.Item(80, 10).Value = 42916
$ws.Cells.Item(80, 11).Value = 42825
$ws.Cells.Item(80, 12).Value = 42735
$ws.Cells.Item(80, 13).Value = 42643

# Row 81
$ws.Cells.Item(81, 4).Value = 56400
$ws.Cells.Item(81, 5).Value = 124200
$ws.Cells.Item(81, 6).Value = 88900
$ws.Cells.Item(81, 7).Value = 42600
$ws.Cells.Item(81, 8).Value = 212200
$ws.Cells.Item(81, 9).Value = 24600
$ws.Cells.Item(81, 10).Value = 62500
$ws.Cells.Item(81, 11).Value = 20900
$ws.Cells.Item(81, 12).Value = 36600
$ws.Cells.Item(81, 13).Value = 32700

# Row 82
$ws.Cells.Item(82, 4).Value = ""
$ws.Cells.Item(82, 5).Value = ""
$ws.Cells.Item(82, 6).Value = ""
$ws.Cells.Item(82, 7).Value = ""
$ws.Cells.Item(82, 8).Value = ""
$ws.Cells.Item(82, 9).Value = ""
$ws.Cells.Item(82, 10).Value = ""
$ws.Cells.Item(82, 11).Value = ""
$ws.Cells.Item(82, 12).Value = ""
$ws.Cells.Item(82, 13).Value = ""

# Row 83
$ws.Cells.Item(83, 4).Value = 28200
$ws.Cells.Item(83, 5).Value = 34000
$ws.Cells.Item(83, 6).Value = 50000
$ws.Cells.Item(83, 7).Value = 25000
$ws.Cells.Item(83, 8).Value = 29000
$ws.Cells.Item(83, 9).Value = 27200
$ws.Cells.Item(83, 10).Value = 57100
$ws.Cells.Item(83, 11).Value = 21500
$ws.Cells.Item(83, 12).Value = 24800
$ws.Cells.Item(83, 13).Value = 22700

# Row 84
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 0

# Row 85
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = 0

# Row 86
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 0

# Row 87
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = 0

# Row 88
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = 0

# Row 89
$ws.Cells.Item(89, 4).Value = 95200
$ws.Cells.Item(89, 5).Value = 121300
$ws.Cells.Item(89, 6).Value = 70500
$ws.Cells.Item(89, 7).Value = 19000
$ws.Cells.Item(89, 8).Value = 47200
$ws.Cells.Item(89, 9).Value = 89300
$ws.Cells.Item(89, 10).Value = 131600
$ws.Cells.Item(89, 11).Value = 91900
$ws.Cells.Item(89, 12).Value = 109800
$ws.Cells.Item(89, 13).Value = 4900

# Row 90
$ws.Cells.Item(90, 4).Value = ""
$ws.Cells.Item(90, 5).Value = ""
$ws.Cells.Item(90, 6).Value = ""
$ws.Cells.Item(90, 7).Value = ""
$ws.Cells.Item(90, 8).Value = ""
$ws.Cells.Item(90, 9).Value = ""
$ws.Cells.Item(90, 10).Value = ""
$ws.Cells.Item(90, 11).Value = ""
$ws.Cells.Item(90, 12).Value = ""
$ws.Cells.Item(90, 13).Value = ""

# Row 91
$ws.Cells.Item(91, 4).Value = -39300
$ws.Cells.Item(91, 5).Value = -18300
$ws.Cells.Item(91, 6).Value = -40500
$ws.Cells.Item(91, 7).Value = -17500
$ws.Cells.Item(91, 8).Value = -16500
$ws.Cells.Item(91, 9).Value = -13900
$ws.Cells.Item(91, 10).Value = -29900
$ws.Cells.Item(91, 11).Value = -15700
$ws.Cells.Item(91, 12).Value = 16100
$ws.Cells.Item(91, 13).Value = -38300

# Row 92
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 0

# Row 93
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 0

# Row 94
$ws.Cells.Item(94, 4).Value = -79700
$ws.Cells.Item(94, 5).Value = -107100
$ws.Cells.Item(94, 6).Value = -44100
$ws.Cells.Item(94, 7).Value = -11500
$ws.Cells.Item(94, 8).Value = -31500
$ws.Cells.Item(94, 9).Value = -75300
$ws.Cells.Item(94, 10).Value = -335100
$ws.Cells.Item(94, 11).Value = -100000
$ws.Cells.Item(94, 12).Value = -26500
$ws.Cells.Item(94, 13).Value = -101200

# Row 95
$ws.Cells.Item(95, 4).Value = ""
$ws.Cells.Item(95, 5).Value = ""
$ws.Cells.Item(95, 6).Value = ""
$ws.Cells.Item(95, 7).Value = ""
$ws.Cells.Item(95, 8).Value = ""
$ws.Cells.Item(95, 9).Value = ""
$ws.Cells.Item(95, 10).Value = ""
$ws.Cells.Item(95, 11).Value = ""
$ws.Cells.Item(95, 12).Value = ""
$ws.Cells.Item(95, 13).Value = ""

# Row 96
$ws.Cells.Item(96, 4).Value = -7100
$ws.Cells.Item(96, 5).Value = -7100
$ws.Cells.Item(96, 6).Value = -14500
$ws.Cells.Item(96, 7).Value = -7300
$ws.Cells.Item(96, 8).Value = -7000
$ws.Cells.Item(96, 9).Value = -7100
$ws.Cells.Item(96, 10).Value = -14200
$ws.Cells.Item(96, 11).Value = -7100
$ws.Cells.Item(96, 12).Value = -6800
$ws.Cells.Item(96, 13).Value = -6800

# Row 97
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 0

# Row 98
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 0

# Row 99
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 0

# Row 100
$ws.Cells.Item(100, 4).Value = -25100
$ws.Cells.Item(100, 5).Value = -40300
$ws.Cells.Item(100, 6).Value = -126900
$ws.Cells.Item(100, 7).Value = -91000
$ws.Cells.Item(100, 8).Value = -25500
$ws.Cells.Item(100, 9).Value = -55400
$ws.Cells.Item(100, 10).Value = -19100
$ws.Cells.Item(100, 11).Value = -9600
$ws.Cells.Item(100, 12).Value = -25800
$ws.Cells.Item(100, 13).Value = 86600

# Row 101
$ws.Cells.Item(101, 4).Value = -2900
$ws.Cells.Item(101, 5).Value = 100
$ws.Cells.Item(101, 6).Value = -4400
$ws.Cells.Item(101, 7).Value = 4200
$ws.Cells.Item(101, 8).Value = 1100
$ws.Cells.Item(101, 9).Value = 2700
$ws.Cells.Item(101, 10).Value = 7000
$ws.Cells.Item(101, 11).Value = 4200
$ws.Cells.Item(101, 12).Value = -7900
$ws.Cells.Item(101, 13).Value = -1300

# Row 102
$ws.Cells.Item(102, 4).Value = -12600
$ws.Cells.Item(102, 5).Value = -26000
$ws.Cells.Item(102, 6).Value = -104800
$ws.Cells.Item(102, 7).Value = -79400
$ws.Cells.Item(102, 8).Value = -8800
$ws.Cells.Item(102, 9).Value = -38700
$ws.Cells.Item(102, 10).Value = -215700
$ws.Cells.Item(102, 11).Value = -13500
$ws.Cells.Item(102, 12).Value = 49600
$ws.Cells.Item(102, 13).Value = -11000
